$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.151.23"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.44%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.306.22"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.39%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.67%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.38%  "

# Row 7
$ws.Range("E7").Value = "  +0.80%  "

# Row 8
$ws.Range("E8").Value = "  +0.01%  "

# Row 9
$ws.Range("E9").Value = "  +1.67%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.68"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.64%  "

# Row 11
$ws.Range("E11").Value = "  +0.25%  "

# Row 12
$ws.Range("E12").Value = "  +0.68%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "17.73"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.96%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.93"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.44%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.665.15"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.24%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.310.57"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.05%  "

# Row 17
$ws.Range("E17").Value = "  -1.15%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.059.84"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.39%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.79"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +10.15%  "

# Row 20
$ws.Range("E20").Value = "  +0.42%  "

# Row 21
$ws.Range("E21").Value = "  +1.82%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.96"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.77%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.17%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.69%  "

# Row 25
$ws.Range("E25").Value = "  +0.47%  "

# Row 26
$ws.Range("E26").Value = "  -0.05%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.33%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "168.03"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.85%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "34.56"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.44%  "

# Row 30
$ws.Range("E30").Value = "  -1.16%  "

# Row 31
$ws.Range("E31").Value = "  +0.50%  "

# Row 32
$ws.Range("E32").Value = "  -0.09%  "

# Row 33
$ws.Range("E33").Value = "  +2.03%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.63"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.93%  "

# Row 35
$ws.Range("E35").Value = "  +0.32%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.42"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.04%  "

# Row 37
$ws.Range("E37").Value = "  -1.27%  "

# Row 38
$ws.Range("E38").Value = "  +1.52%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.79"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.25%  "

# Row 40
$ws.Range("E40").Value = "  -0.17%  "

# Row 41
$ws.Range("E41").Value = "  +0.80%  "

# Row 42
$ws.Range("E42").Value = "  +3.26%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.982.68"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.42%  "

# Row 44
$ws.Range("E44").Value = "  -4.23%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.28"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.53%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.85"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.31%  "

# Row 47
$ws.Range("E47").Value = "  +2.47%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "55.39"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.21%  "

# Row 49
$ws.Range("E49").Value = "  +4.50%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.531.49"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.21%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "70.84"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.29%  "
